$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "Test for Git"
$ws.Range("B24:D24").Style = "Normal"
